$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.7263964236600271
$ws.Range("C2").Value = 0.7334881754064316
$ws.Range("D2").Value = 1.472823318998768
$ws.Range("E2").Value = 0.7483408004336237
$ws.Range("F2").Value = 0.740830947129459
$ws.Range("G2").Value = 0.7653563574070629
$ws.Range("H2").Value = 0.7487255867942739
$ws.Range("B3").Value = 0.7497963153137237
$ws.Range("C3").Value = 0.7369496738843104
$ws.Range("D3").Value = 1.409469658895177
$ws.Range("E3").Value = 0.7522648609778715
$ws.Range("F3").Value = 0.7557067841118835
$ws.Range("G3").Value = 0.767545930173578
$ws.Range("H3").Value = 0.7524391359638646
$ws.Range("B4").Value = 0.6784402355552254
$ws.Range("C4").Value = 0.7602560484723868
$ws.Range("D4").Value = 1.082814381615641
$ws.Range("E4").Value = 0.7427419605067365
$ws.Range("F4").Value = 0.695716271214237
$ws.Range("G4").Value = 0.7489833093671484
$ws.Range("H4").Value = 0.7427761591611421
$ws.Range("B5").Value = 0.7845089806605176
$ws.Range("C5").Value = 0.7861785817351283
$ws.Range("D5").Value = 0.3980077702589124
$ws.Range("E5").Value = 0.7596649851233738
$ws.Range("F5").Value = 0.7775079597791049
$ws.Range("G5").Value = 0.7722547943492591
$ws.Range("H5").Value = 0.7593728668952565
$ws.Range("B6").Value = 0.7756754437877884
$ws.Range("C6").Value = 0.7846734571989462
$ws.Range("D6").Value = 0.2390097767754088
$ws.Range("E6").Value = 0.7583094459532475
$ws.Range("F6").Value = 0.7695806858902364
$ws.Range("G6").Value = 0.766576297593929
$ws.Range("H6").Value = 0.7579578215545946
$ws.Range("B7").Value = 0.7937543922314764
$ws.Range("C7").Value = 0.7892787435559668
$ws.Range("D7").Value = 0.3706086512129885
$ws.Range("E7").Value = 0.7601619967263648
$ws.Range("F7").Value = 0.7866323612145677
$ws.Range("G7").Value = 0.7730594823287936
$ws.Range("H7").Value = 0.7598586426095999
$ws.Range("B8").Value = 0.7384478743081042
$ws.Range("C8").Value = 0.709022198366438
$ws.Range("D8").Value = 0.3085576887316853
$ws.Range("E8").Value = 0.7391238751082457
$ws.Range("F8").Value = 0.7321063414733795
$ws.Range("G8").Value = 0.7053807860549471
$ws.Range("H8").Value = 0.7388605506954972
$ws.Range("B9").Value = 0.8003621903725986
$ws.Range("C9").Value = 0.7651098808327893
$ws.Range("D9").Value = 0.473483256811177
$ws.Range("E9").Value = 0.7585389402505501
$ws.Range("F9").Value = 0.793596799859918
$ws.Range("G9").Value = 0.7539730207653869
$ws.Range("H9").Value = 0.7583282497459432
$ws.Range("B10").Value = 0.7818163400102417
$ws.Range("C10").Value = 0.7863713395062196
$ws.Range("D10").Value = 1.275327798525587
$ws.Range("E10").Value = 0.7585068958190951
$ws.Range("F10").Value = 0.7849380575792405
$ws.Range("G10").Value = 0.790092921855691
$ws.Range("H10").Value = 0.758778634037358
$ws.Range("B11").Value = 0.7742004243341288
$ws.Range("C11").Value = 0.78756876434718
$ws.Range("D11").Value = 1.162935609077621
$ws.Range("E11").Value = 0.7591333756474891
$ws.Range("F11").Value = 0.7772108670126719
$ws.Range("G11").Value = 0.7882143195690925
$ws.Range("H11").Value = 0.7592666848883107
$ws.Range("B12").Value = 0.6573517222229365
$ws.Range("C12").Value = 0.7885821271855722
$ws.Range("D12").Value = 0.614786416514029
$ws.Range("E12").Value = 0.7528757606733335
$ws.Range("F12").Value = 0.6700888040979659
$ws.Range("G12").Value = 0.7731175937837654
$ws.Range("H12").Value = 0.7525739072257516
$ws.Range("B13").Value = 0.77706033667528
$ws.Range("C13").Value = 0.7855739167611877
$ws.Range("D13").Value = 1.115428100200711
$ws.Range("E13").Value = 0.7577073778010187
$ws.Range("F13").Value = 0.7789775295704988
$ws.Range("G13").Value = 0.7854760679862451
$ws.Range("H13").Value = 0.7578752676462324
